$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 34894.363
$ws.Range("I69").Value = 4912
$ws.Range("K69").Value = 14736
$ws.Range("M69").Value = -13862
$ws.Range("H72").Value = 34894.363
$ws.Range("I72").Value = 4912
$ws.Range("K72").Value = 44208
$ws.Range("M72").Value = -39840
$ws.Range("H76").Value = 6250.6665
$ws.Range("I76").Value = 6250.6665
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 6250.6665
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -5935.6665
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 6250.6665
$ws.Range("I79").Value = 6250.6665
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 6250.6665
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -5158.6665
$ws.Range("N79").ClearContents()
$ws.Range("H100").Value = 3856.1304
$ws.Range("I100").Value = 1542.625
$ws.Range("J100").Value = 9144.143
$ws.Range("K100").Value = 1542.625
$ws.Range("L100").Value = 9144.143
$ws.Range("M100").Value = -1001.625
$ws.Range("N100").Value = -10226.143
$ws.Range("H132").Value = 955.45654
$ws.Range("J132").Value = 1188.6
$ws.Range("L132").Value = 3565.8
$ws.Range("N132").Value = -8625.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1845046.5
$ws.Range("I32").Value = 9731.861999999999
$ws.Range("J32").Value = 15192789
$ws.Range("K32").Value = 9731.861999999999
$ws.Range("L32").Value = 15192789
$ws.Range("M32").Value = -9444.861999999999
$ws.Range("N32").Value = -15193363
$ws.Range("H132").Value = 8972.846
$ws.Range("I132").Value = 8960.777
$ws.Range("K132").Value = 26882.331
$ws.Range("M132").Value = -24352.331

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4114.846
$ws.Range("I105").Value = 3449.75
$ws.Range("J105").Value = 4410.4443
$ws.Range("K105").Value = 3449.75
$ws.Range("L105").Value = 4410.4443
$ws.Range("M105").Value = -1702.75
$ws.Range("N105").Value = -7904.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2063.6667
$ws.Range("I2").Value = 1681
$ws.Range("J2").Value = 2255
$ws.Range("K2").Value = 1681
$ws.Range("L2").Value = 2255
$ws.Range("M2").Value = -1568
$ws.Range("N2").Value = -2481
$ws.Range("H59").Value = 999999
$ws.Range("I59").Value = 999999
$ws.Range("K59").Value = 999999
$ws.Range("M59").Value = -998854
$ws.Range("H105").Value = 956.44446
$ws.Range("I105").Value = 902.7143
$ws.Range("K105").Value = 902.7143
$ws.Range("M105").Value = 844.2857
$ws.Range("H132").Value = 4359.684
$ws.Range("I132").Value = 3521.6667
$ws.Range("K132").Value = 10565.0001
$ws.Range("M132").Value = -8035.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5832.95
$ws.Range("I5").Value = 575.3125
$ws.Range("J5").Value = 26863.5
$ws.Range("K5").Value = 1725.9375
$ws.Range("L5").Value = 80590.5
$ws.Range("M5").Value = -1613.9375
$ws.Range("N5").Value = -80814.5
$ws.Range("H37").Value = 239035.58
$ws.Range("J37").Value = 239035.58
$ws.Range("L37").Value = 717106.74
$ws.Range("N37").Value = -717330.74
$ws.Range("H92").Value = 549.6667
$ws.Range("I92").Value = 299
$ws.Range("K92").Value = 897
$ws.Range("M92").Value = 351
$ws.Range("H131").Value = 4492.9443
$ws.Range("I131").Value = 3911.2856
$ws.Range("J131").Value = 4863.091
$ws.Range("K131").Value = 11733.8568
$ws.Range("L131").Value = 14589.273
$ws.Range("M131").Value = -6693.856800000001
$ws.Range("N131").Value = -24669.273
$ws.Range("H135").Value = 5832.95
$ws.Range("I135").Value = 575.3125
$ws.Range("J135").Value = 26863.5
$ws.Range("K135").Value = 5177.8125
$ws.Range("L135").Value = 241771.5
$ws.Range("M135").Value = -2642.8125
$ws.Range("N135").Value = -246841.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 8065168
$ws.Range("I97").Value = 642.2778
$ws.Range("K97").Value = 642.2778
$ws.Range("M97").Value = -146.2778
$ws.Range("H122").Value = 4233.4
$ws.Range("I122").Value = 1950.5
$ws.Range("J122").Value = 5755.3335
$ws.Range("K122").Value = 5851.5
$ws.Range("L122").Value = 17266.0005
$ws.Range("M122").Value = -3401.5
$ws.Range("N122").Value = -22166.0005
$ws.Range("H126").Value = 2516.8333
$ws.Range("J126").Value = 4501.3
$ws.Range("L126").Value = 13503.9
$ws.Range("N126").Value = -18443.9
$ws.Range("H132").Value = 5230.55
$ws.Range("I132").Value = 4663.517
$ws.Range("K132").Value = 13990.551
$ws.Range("M132").Value = -11460.551

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 19499.334
$ws.Range("J97").Value = 19499.334
$ws.Range("L97").Value = 19499.334
$ws.Range("N97").Value = -21481.334
$ws.Range("H100").Value = 89050.53999999999
$ws.Range("I100").Value = 125295.22
$ws.Range("K100").Value = 125295.22
$ws.Range("M100").Value = -124754.22
$ws.Range("H108").Value = 374999.34
$ws.Range("J108").Value = 374999.34
$ws.Range("L108").Value = 374999.34
$ws.Range("N108").Value = -382679.34
$ws.Range("H132").Value = 4855.2163
$ws.Range("I132").Value = 4317.048
$ws.Range("J132").Value = 5561.5625
$ws.Range("K132").Value = 12951.144
$ws.Range("L132").Value = 16684.6875
$ws.Range("M132").Value = -10421.144
$ws.Range("N132").Value = -21744.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 4910
$ws.Range("I29").Value = 4910
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 4910
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -4620
$ws.Range("N29").ClearContents()
$ws.Range("H96").Value = 64588.625
$ws.Range("I96").Value = 126701.25
$ws.Range("K96").Value = 126701.25
$ws.Range("M96").Value = -125328.25
$ws.Range("H132").Value = 2489.4375
$ws.Range("I132").Value = 1427.5
$ws.Range("K132").Value = 4282.5
$ws.Range("M132").Value = -1752.5
$ws.Range("H136").Value = 10002897
$ws.Range("I136").Value = 14708199
$ws.Range("J136").Value = 4130.9375
$ws.Range("K136").Value = 44124597
$ws.Range("L136").Value = 12392.8125
$ws.Range("M136").Value = -44122047
$ws.Range("N136").Value = -17492.8125
